$wb = $excel.ActiveWorkbook

# row 12 in sheet ALC (diff @ line 1238)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 385.125
$ws.Range("I12").Value = 240.25
$ws.Range("J12").Value = 530
$ws.Range("K12").Value = 240.25
$ws.Range("L12").Value = 530
$ws.Range("M12").Value = -70.25
$ws.Range("N12").Value = -870

# row 17 in sheet ALC (diff @ line 1489)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1803515.9
$ws.Range("J17").Value = 1803515.9
$ws.Range("L17").Value = 5410547.699999999
$ws.Range("N17").Value = -5410883.699999999

# row 100 in sheet ALC (diff @ line 5667)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1561
$ws.Range("I100").Value = 1668.3334
$ws.Range("J100").Value = 1400
$ws.Range("K100").Value = 1668.3334
$ws.Range("L100").Value = 1400
$ws.Range("M100").Value = -1127.3334
$ws.Range("N100").Value = -2482

# row 129 in sheet ALC (diff @ line 7124)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1091.4595
$ws.Range("I129").Value = 491.42856
$ws.Range("K129").Value = 1474.28568
$ws.Range("M129").Value = 3525.71432

# row 137 in sheet ALC (diff @ line 7528)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1007.2222
$ws.Range("I137").Value = 999.8
$ws.Range("J137").Value = 1100
$ws.Range("K137").Value = 2999.4
$ws.Range("L137").Value = 3300
$ws.Range("M137").Value = -449.3999999999996
$ws.Range("N137").Value = -8400

# row 138 in sheet ALC (diff @ line 7580)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2068.015
$ws.Range("I138").Value = 1452.1111
$ws.Range("J138").Value = 2294.2654
$ws.Range("K138").Value = 4356.3333
$ws.Range("L138").Value = 6882.796200000001
$ws.Range("M138").Value = 783.6666999999998
$ws.Range("N138").Value = -17162.7962

# row 102 in sheet ARM (diff @ line 12776)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1732.25
$ws.Range("I102").Value = 1758.2858
$ws.Range("K102").Value = 1758.2858
$ws.Range("M102").Value = -136.2858000000001

# row 132 in sheet ARM (diff @ line 14243)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3084.25
$ws.Range("I132").Value = 1892.8695
$ws.Range("J132").Value = 5192.077
$ws.Range("K132").Value = 5678.6085
$ws.Range("L132").Value = 15576.231
$ws.Range("M132").Value = -3148.6085
$ws.Range("N132").Value = -20636.231

# row 61 in sheet BSM (diff @ line 17721)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 20000
$ws.Range("J61").Value = 20000
$ws.Range("L61").Value = 20000
$ws.Range("N61").Value = -20626

# row 105 in sheet BSM (diff @ line 19889)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 17440.5
$ws.Range("I105").Value = 23067.7
$ws.Range("J105").Value = 3372.5
$ws.Range("K105").Value = 23067.7
$ws.Range("L105").Value = 3372.5
$ws.Range("M105").Value = -21320.7
$ws.Range("N105").Value = -6866.5

# row 31 in sheet CRP (diff @ line 23226)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3751
$ws.Range("I31").Value = 2704.6667
$ws.Range("J31").Value = 4498.381
$ws.Range("K31").Value = 2704.6667
$ws.Range("L31").Value = 4498.381
$ws.Range("M31").Value = -2409.6667
$ws.Range("N31").Value = -5088.381

# row 34 in sheet CRP (diff @ line 23379)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3751
$ws.Range("I34").Value = 2704.6667
$ws.Range("J34").Value = 4498.381
$ws.Range("K34").Value = 2704.6667
$ws.Range("L34").Value = 4498.381
$ws.Range("M34").Value = -2502.6667
$ws.Range("N34").Value = -4902.381

# row 121 in sheet CUL (diff @ line 34785)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 8864.65
$ws.Range("I121").Value = 610
$ws.Range("J121").Value = 11616.2
$ws.Range("K121").Value = 1830
$ws.Range("L121").Value = 34848.60000000001
$ws.Range("M121").Value = -520
$ws.Range("N121").Value = -37468.60000000001

# row 131 in sheet CUL (diff @ line 35299)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2440041.2
$ws.Range("J131").Value = 1153.8788
$ws.Range("L131").Value = 3461.6364
$ws.Range("N131").Value = -13541.6364

# row 132 in sheet CUL (diff @ line 35351)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2708.5833
$ws.Range("I132").Value = 3613.75
$ws.Range("J132").Value = 2256
$ws.Range("K132").Value = 32523.75
$ws.Range("L132").Value = 20304
$ws.Range("M132").Value = -29993.75
$ws.Range("N132").Value = -25364

# row 126 in sheet GSM (diff @ line 42038)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4620.2085
$ws.Range("I126").Value = 7953.375
$ws.Range("J126").Value = 2953.625
$ws.Range("K126").Value = 23860.125
$ws.Range("L126").Value = 8860.875
$ws.Range("M126").Value = -21390.125
$ws.Range("N126").Value = -13800.875

# row 132 in sheet GSM (diff @ line 42326)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2984.7222
$ws.Range("I132").Value = 3383.889
$ws.Range("K132").Value = 10151.667
$ws.Range("M132").Value = -7621.667000000001

# row 7 in sheet LTW (diff @ line 43155)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 202280.8
$ws.Range("I7").Value = 202280.8
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 202280.8
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -202168.8
$ws.Range("N7").ClearContents()

# row 40 in sheet LTW (diff @ line 44781)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2787.8125
$ws.Range("I40").Value = 2935.7144
$ws.Range("J40").Value = 1752.5
$ws.Range("K40").Value = 2935.7144
$ws.Range("L40").Value = 1752.5
$ws.Range("M40").Value = -2799.7144
$ws.Range("N40").Value = -2024.5

# row 46 in sheet LTW (diff @ line 45072)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1306.5
$ws.Range("I46").Value = 1075
$ws.Range("J46").Value = 1422.25
$ws.Range("K46").Value = 1075
$ws.Range("L46").Value = 1422.25
$ws.Range("M46").Value = -887
$ws.Range("N46").Value = -1798.25

# row 82 in sheet LTW (diff @ line 46833)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 929128.4399999999
$ws.Range("I82").Value = 1667791.5
$ws.Range("K82").Value = 1667791.5
$ws.Range("M82").Value = -1667430.5

# row 85 in sheet LTW (diff @ line 46983)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 929128.4399999999
$ws.Range("I85").Value = 1667791.5
$ws.Range("K85").Value = 1667791.5
$ws.Range("M85").Value = -1666543.5

# row 100 in sheet LTW (diff @ line 47709)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1264.2106
$ws.Range("I100").Value = 1001.53845
$ws.Range("K100").Value = 1001.53845
$ws.Range("M100").Value = -460.53845

# row 122 in sheet LTW (diff @ line 48766)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 592176.3
$ws.Range("I122").Value = 4161.3076
$ws.Range("J122").Value = 2503225
$ws.Range("K122").Value = 12483.9228
$ws.Range("L122").Value = 7509675
$ws.Range("M122").Value = -10033.9228
$ws.Range("N122").Value = -7514575

# row 126 in sheet LTW (diff @ line 48959)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 202280.8
$ws.Range("I126").Value = 202280.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 606842.3999999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -604372.3999999999
$ws.Range("N126").ClearContents()

# row 132 in sheet LTW (diff @ line 49247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9123.647000000001
$ws.Range("I132").Value = 10907.923
$ws.Range("J132").Value = 3324.75
$ws.Range("K132").Value = 32723.769
$ws.Range("L132").Value = 9974.25
$ws.Range("M132").Value = -30193.769
$ws.Range("N132").Value = -15034.25

# row 96 in sheet WVR (diff @ line 54410)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 934
$ws.Range("I96").Value = 960.0952
$ws.Range("J96").Value = 879.2
$ws.Range("K96").Value = 960.0952
$ws.Range("L96").Value = 879.2
$ws.Range("M96").Value = 412.9048
$ws.Range("N96").Value = -3625.2
